$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 2335
    $ws.Range("F3").Value = 1793
    $ws.Range("F6").Value = 1001
}
